# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Septiembre de 2020 a las 21:01"

# --- Refresh case numbers for countries whose row position does not move ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 6448652
$ws.Range("C4").Value = 17500
$ws.Range("D4").Value = 3714315
$ws.Range("E4").Value = 2541296
$ws.Range("G4").Value = 223
$ws.Range("H4").Value = 193041

# India (row 5)
$ws.Range("B5").Value = 4202562
$ws.Range("C5").Value = 91723
$ws.Range("D5").Value = 3247297
$ws.Range("E5").Value = 883526

# Irak (row 23)
$ws.Range("B23").Value = 260370
$ws.Range("C23").Value = 3651
$ws.Range("D23").Value = 198560
$ws.Range("E23").Value = 54298
$ws.Range("G23").Value = 90
$ws.Range("H23").Value = 7512

# Alemania (row 24)
$ws.Range("B24").Value = 251456
$ws.Range("C24").Value = 400
$ws.Range("E24").Value = 15847

# Suazilandia (row 113)
$ws.Range("B113").Value = 4853
$ws.Range("C113").Value = 34
$ws.Range("D113").Value = 3951
$ws.Range("E113").Value = 808

# Islandia (row 145)
$ws.Range("B145").Value = 2141
$ws.Range("C145").Value = 5
$ws.Range("D145").Value = 2050
$ws.Range("E145").Value = 81

# Yemen (row 149)
$ws.Range("B149").Value = 1987
$ws.Range("C149").Value = 4
$ws.Range("D149").Value = 1200
$ws.Range("E149").Value = 215

# Republica de Chipre (row 153)
$ws.Range("B153").Value = 1509
$ws.Range("C153").Value = 2
$ws.Range("E153").Value = 251

# --- Namibia now reported ahead of Gabon / Maldivas (rows 99-101 keep their
#     position, but the country order - and therefore the figures shown in
#     each row - shifts down by one) ---

$ws.Range("A99").Value = "Namibia"
$ws.Range("B99").Value = 8685
$ws.Range("C99").Value = 171
$ws.Range("D99").Value = 3786
$ws.Range("E99").Value = 4810
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 89

$ws.Range("A100").Value = "Gabon"
$ws.Range("B100").Value = 8601
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 7424
$ws.Range("E100").Value = 1124
$ws.Range("H100").Value = 53

$ws.Range("A101").Value = "Maldivas"
$ws.Range("B101").Value = 8584
$ws.Range("C101").Value = 98
$ws.Range("D101").Value = 5936
$ws.Range("E101").Value = 2619
$ws.Range("H101").Value = 29

# --- Siria now reported ahead of Eslovenia / Sri Lanka (rows 129-131 keep
#     their position, figures shift down by one row) ---

$ws.Range("A129").Value = "Siria"
$ws.Range("B129").Value = 3171
$ws.Range("C129").Value = 67
$ws.Range("D129").Value = 730
$ws.Range("E129").Value = 2307
$ws.Range("G129").Value = 4
$ws.Range("H129").Value = 134

$ws.Range("A130").Value = "Eslovenia"
$ws.Range("B130").Value = 3165
$ws.Range("C130").Value = 43
$ws.Range("D130").Value = 2483
$ws.Range("E130").Value = 547
$ws.Range("H130").Value = 135

$ws.Range("A131").Value = "Sri Lanka"
$ws.Range("B131").Value = 3122
$ws.Range("C131").Value = 1
$ws.Range("D131").Value = 2925
$ws.Range("E131").Value = 185
$ws.Range("H131").Value = 12
